# "mend the equip effect" - EquipAddon.xlsx
# - Fix the Format strings for existing "attack"/"life" equip-addon rows so
#   they describe tower attack/life bonuses instead of generic ones.
# - Add two new equip-addon rows (attack speed / range).
# - Resize the worksheet Table to cover the newly added rows.
# - Update the active selection and the workbook's background theme color.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the table (ListObject) so the two new rows become part of it.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E7"))

# Row 4 ("攻击" / attack): fix the Format text.
$ws.Range("C4").Value2 = "防御塔攻击+{0}%"

# Row 5 ("生命" / life): fix the Format text.
$ws.Range("C5").Value2 = "防御塔生命+{0}%"

# New rows: set the Name (column B) cells before the Format (column C)
# cells so the shared-string table is built in the same order as the
# reference workbook.
$ws.Range("B6").Value2 = "射速"
$ws.Range("B7").Value2 = "射程"
$ws.Range("C6").Value2 = "防御塔射速+{0}"
$ws.Range("C7").Value2 = "防御塔射程+{0}"

$ws.Range("A6").Value2 = 3
$ws.Range("D6").Value2 = "attr"
$ws.Range("E6").Value2 = 2

$ws.Range("A7").Value2 = 4
$ws.Range("D7").Value2 = "attr"
$ws.Range("E7").Value2 = 2

# Match the saved selection in the worksheet.
[void]$ws.Range("C7").Select()

# Update the workbook theme's "Background 1" / window color.
$theme = $wb.Theme
$cs = $theme.ThemeColorScheme
$cs.Colors(2).RGB = 13494986

Write-Output "done"
